# Add new columns I (I0) and J (IF) to the sheet, matching the styling
# of the existing header row and filling in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, centered/top alignment)
# from the existing H1 header cell onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (2-81): I0 and IF values ---
$data = @(
    @(2, 5, 6),
    @(3, 7, 7),
    @(4, 7, 7),
    @(5, 8, 8),
    @(6, 7, 7),
    @(7, 8, 8),
    @(8, 6, 6),
    @(9, 5, 6),
    @(10, 8, 8),
    @(11, 7, 7),
    @(12, 8, 8),
    @(13, 7, 7),
    @(14, 5, 5),
    @(15, 6, 6),
    @(16, 9, 10),
    @(17, 6, 6),
    @(18, 8, 8),
    @(19, 5, 5),
    @(20, 7, 7),
    @(21, 7, 7),
    @(22, 7, 7),
    @(23, 8, 8),
    @(24, 8, 8),
    @(25, 7, 7),
    @(26, 7, 7),
    @(27, 10, 10),
    @(28, 9, 9),
    @(29, 9, 9),
    @(30, 7, 8),
    @(31, 7, 7),
    @(32, 8, 8),
    @(33, 7, 7),
    @(34, 8, 8),
    @(35, 9, 9),
    @(36, 9, 9),
    @(37, 8, 8),
    @(38, 11, 11),
    @(39, 6, 6),
    @(40, 8, 8),
    @(41, 9, 9),
    @(42, 9, 9),
    @(43, 7, 7),
    @(44, 6, 6),
    @(45, 6, 7),
    @(46, 5, 6),
    @(47, 8, 8),
    @(48, 7, 7),
    @(49, 6, 6),
    @(50, 9, 9),
    @(51, 8, 8),
    @(52, 8, 8),
    @(53, 11, 11),
    @(54, 9, 10),
    @(55, 7, 7),
    @(56, 8, 8),
    @(57, 9, 9),
    @(58, 6, 7),
    @(59, 7, 7),
    @(60, 7, 8),
    @(61, 8, 8),
    @(62, 6, 7),
    @(63, 8, 9),
    @(64, 7, 7),
    @(65, 8, 8),
    @(66, 8, 8),
    @(67, 9, 9),
    @(68, 7, 8),
    @(69, 7, 7),
    @(70, 8, 8),
    @(71, 7, 7),
    @(72, 7, 7),
    @(73, 8, 8),
    @(74, 6, 7),
    @(75, 8, 8),
    @(76, 8, 8),
    @(77, 9, 9),
    @(78, 7, 7),
    @(79, 8, 8),
    @(80, 3, 3),
    @(81, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
